$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet tracks, per shop (row), a water-bottle delivery cycle:
#   D = total days in the cycle
#   E = days remaining (= D - elapsed days since F)
#   F = cycle start date, encoded as an 8-digit integer yyyyMMdd
#
# This edit advances the reference "today" by one day (from 2025-12-10 to
# 2025-12-11), so every row's remaining count (E) drops by one. If a row's
# remaining count would hit zero or below (cycle finished), the row is
# restocked instead: the cycle restarts today, i.e. F becomes the new
# "today" (20251211) and E resets back up to the row's own total D.

function DateToDayNumber($y, $m, $d) {
    if ($m -le 2) {
        $y = $y - 1
        $m = $m + 12
    }
    $era = [math]::Floor($y / 400)
    $yoe = $y - $era * 400
    $doy = [math]::Floor((153 * ($m - 3) + 2) / 5) + $d - 1
    $doe = $yoe * 365 + [math]::Floor($yoe / 4) - [math]::Floor($yoe / 100) + $doy
    return $era * 146097 + $doe
}

$todayAfterNum = DateToDayNumber 2025 12 11
$todayAfterCode = 20251211

for ($r = 2; $r -le 99; $r++) {
    $d = $ws.Cells.Item($r, 4).Value()
    $e = $ws.Cells.Item($r, 5).Value()
    $f = $ws.Cells.Item($r, 6).Value()

    if ($d -eq $null -or $e -eq $null -or $f -eq $null) {
        continue
    }

    $fStr = [string]([int]$f)
    if ($fStr.Length -ne 8) {
        # Malformed date code (e.g. a stray digit) - leave this row untouched.
        continue
    }

    $year = [int]$fStr.Substring(0, 4)
    $month = [int]$fStr.Substring(4, 2)
    $day = [int]$fStr.Substring(6, 2)

    if ($month -lt 1 -or $month -gt 12 -or $day -lt 1 -or $day -gt 31) {
        continue
    }

    $startNum = DateToDayNumber $year $month $day
    $elapsedAfter = $todayAfterNum - $startNum
    $newE = $d - $elapsedAfter

    if ($newE -le 0) {
        # Cycle finished - restock: cycle restarts as of the new "today".
        $ws.Cells.Item($r, 5).Value = $d
        $ws.Cells.Item($r, 6).Value = $todayAfterCode
    } else {
        $ws.Cells.Item($r, 5).Value = $newE
    }
}
